# Finanzbericht_Vorlage.xlsx - "xlxs reader funktional (ungetestet)"
#
# The reader script apparently produced three rows that don't belong in the
# Bilanz/Erfolgsrechnung template:
#   - "Wertschriften"              (old row 10, under "Flüssige Mittel")
#   - "Bankshuld"                  (old row 23, under "Fremdkapital")
#   - "Verluste aus Forderungen"   (old row 48, under "Ertrag")
#
# Removing the three rows lets Excel shift everything below them up and
# automatically repair the SUM(...) ranges / totals that reference them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete bottom-to-top so the row numbers used here stay valid for each
# subsequent call (deleting a lower row does not renumber the ones above it).
$ws.Rows(48).Delete()   # "Verluste aus Forderungen" (-34)
$ws.Rows(23).Delete()   # "Bankshuld" (80)
$ws.Rows(10).Delete()   # "Wertschriften" (44)

# Re-assign the "Total Aufwand" side totals (C/D) as one multi-cell formula so
# Excel keeps storing them as a shared formula group, same as the original
# workbook did before the row shift.
$ws.Range("C40:D40").Formula = "=SUM(C33:C39)"

# Restore cursor / scroll position to roughly where the author left it
# (top of the Fremdkapital/Eigenkapital block, with the active cell sitting
# further down in the Erfolgsrechnung section).
[void]$ws.Range("A17").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1

[void]$ws.Range("J43").Select()
